$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.612.04'
$ws.Range("E2").Value = '  -0.69%  '

$ws.Range("D3").Value = '2.284.38'
$ws.Range("E3").Value = '  -2.56%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.19%  '

$ws.Range("E7").Value = '  -1.51%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  -3.13%  '

$ws.Range("E10").Value = '  -5.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.25%  '

$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '2.640.17'
$ws.Range("E15").Value = '  -2.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("D17").Value = '2.324.16'
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("E18").Value = '  -1.61%  '

$ws.Range("D19").Value = '42.538.37'
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '

$ws.Range("E25").Value = '  -3.38%  '

$ws.Range("E26").Value = '  -3.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.74%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '165.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.18%  '

$ws.Range("E32").Value = '  -3.18%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -3.70%  '

$ws.Range("E35").Value = '  -3.94%  '

$ws.Range("E36").Value = '  -4.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.67%  '

$ws.Range("E38").Value = '  -6.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.05%  '

$ws.Range("E40").Value = '  -6.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0992'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.23%  '

$ws.Range("E42").Value = '  -2.40%  '

$ws.Range("E43").Value = '  -1.95%  '

$ws.Range("D44").Value = '1.961.61'
$ws.Range("E44").Value = '  -2.96%  '

$ws.Range("E45").Value = '  -1.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.74%  '

$ws.Range("E48").Value = '  -6.20%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.51%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.509.82'
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.58%  '
